# Update cryptocurrency price/volume figures (GitHub Actions refresh).
# Each cell is forced to Text via NumberFormat "@" before the write so
# numeric-looking strings (e.g. "520.60", "1.00") are not silently
# coerced to the Number type by Excel's smart-entry parser, then
# ClearFormats() strips the temporary "@" style back off so no stray
# cell-format/style is left behind (matches the source diff, which is
# a pure text-content change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}


Set-TextValue 'D2' '59.089.02'
Set-TextValue 'E2' '  +1.36%  '
Set-TextValue 'D3' '2.625.79'
Set-TextValue 'E3' '  +3.36%  '
Set-TextValue 'E4' '  +0.21%  '
Set-TextValue 'D5' '520.60'
Set-TextValue 'E5' '  +3.06%  '
Set-TextValue 'D6' '145.90'
Set-TextValue 'D7' '0.998'
Set-TextValue 'E7' '  -0.03%  '
Set-TextValue 'D8' '0.568'
Set-TextValue 'E8' '  +0.62%  '
Set-TextValue 'D9' '2.648.74'
Set-TextValue 'E9' '  +4.06%  '
Set-TextValue 'E10' '  +3.62%  '
Set-TextValue 'E11' '  +3.82%  '
Set-TextValue 'E12' '  +2.35%  '
Set-TextValue 'E13' '  -1.31%  '
Set-TextValue 'D14' '3.109.68'
Set-TextValue 'E14' '  +4.09%  '
Set-TextValue 'D15' '59.077.39'
Set-TextValue 'E15' '  +1.37%  '
Set-TextValue 'D16' '21.05'
Set-TextValue 'E16' '  +1.97%  '
Set-TextValue 'E17' '  +2.03%  '
Set-TextValue 'D18' '2.641.27'
Set-TextValue 'E18' '  +3.50%  '
Set-TextValue 'D19' '349.40'
Set-TextValue 'E19' '  +2.48%  '
Set-TextValue 'D20' '4.52'
Set-TextValue 'E20' '  +0.43%  '
Set-TextValue 'D21' '10.35'
Set-TextValue 'E21' '  +3.18%  '
Set-TextValue 'D22' '6.18'
Set-TextValue 'E22' '  +4.34%  '
Set-TextValue 'D23' '1.00'
Set-TextValue 'E23' '  +0.20%  '
Set-TextValue 'D24' '61.71'
Set-TextValue 'E24' '  +1.82%  '
Set-TextValue 'E25' '  +2.58%  '
Set-TextValue 'E26' '  +3.44%  '
Set-TextValue 'E27' '  -0.39%  '
Set-TextValue 'D28' '0.0₃0809'
Set-TextValue 'E28' '  +3.52%  '
Set-TextValue 'D29' '7.13'
Set-TextValue 'E29' '  +3.12%  '
Set-TextValue 'D30' '1.00'
Set-TextValue 'E30' '  +0.03%  '
Set-TextValue 'D31' '6.26'
Set-TextValue 'E31' '  +7.99%  '
Set-TextValue 'D32' '19.01'
Set-TextValue 'E32' '  +3.01%  '
Set-TextValue 'E33' '  +3.60%  '
Set-TextValue 'D34' '150.15'
Set-TextValue 'E34' '  +1.03%  '
Set-TextValue 'D35' '0.977'
Set-TextValue 'E35' '  +8.76%  '
Set-TextValue 'D36' '4.02'
Set-TextValue 'E36' '  +3.86%  '
Set-TextValue 'E37' '  +3.13%  '
Set-TextValue 'D38' '36.75'
Set-TextValue 'E38' '  +2.51%  '
Set-TextValue 'D39' '0.847'
Set-TextValue 'E39' '  +3.94%  '
Set-TextValue 'D40' '3.70'
Set-TextValue 'E40' '  +5.37%  '
Set-TextValue 'E41' '  +2.66%  '
Set-TextValue 'D42' '277.67'
Set-TextValue 'E42' '  -0.56%  '
Set-TextValue 'D43' '0.994'
Set-TextValue 'E43' '  -0.34%  '
Set-TextValue 'E44' '  -0.54%  '
Set-TextValue 'D45' '0.607'
Set-TextValue 'E45' '  +1.41%  '
Set-TextValue 'D46' '19.68'
Set-TextValue 'E46' '  +6.03%  '
Set-TextValue 'D47' '0.0524'
Set-TextValue 'E47' '  -0.94%  '
Set-TextValue 'E48' '  +2.20%  '
Set-TextValue 'E49' '  +0.05%  '
Set-TextValue 'D50' '1.984.88'
Set-TextValue 'E50' '  +4.96%  '
Set-TextValue 'E51' '  +3.18%  '
